$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The match previously stored in row 7 (Napredak vs Zeleznicar Pancevo) and
#    the one in row 8 (Crvena zvezda vs Vojvodina) were swapped - row 7 now
#    holds the Crvena zvezda game and row 8 the Napredak game. Columns A:E
#    (index / pais / torneio / temporada / data_partida) stay as-is since
#    both games share the same date; only F:V (home..url_partida) swap.
$row7vals = $ws.Range("F7:V7").Value()
$row8vals = $ws.Range("F8:V8").Value()

$ws.Range("F7:V7").Value = $row8vals
$ws.Range("F8:V8").Value = $row7vals

# 2) A new match row was appended at the bottom (row 112). Copy the
#    formatting of the last existing row (111) down so the new row keeps the
#    same styles (bold/bordered index cell, date-formatted data_partida
#    cell, etc.), then fill in the new row's values.
$ws.Range("A111:V111").Copy()
$ws.Range("A112:V112").PasteSpecial(-4122)

$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = "serbia"
$ws.Cells.Item(112, 3).Value = "super-liga"
$ws.Cells.Item(112, 4).Value = "2023-2024"
$ws.Cells.Item(112, 5).Value = 45254.70833333334
$ws.Cells.Item(112, 6).Value = "IMT Novi Beograd"
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = "Sp. Subotica"
$ws.Cells.Item(112, 9).Value = 1
$ws.Cells.Item(112, 10).Value = 1.83
$ws.Cells.Item(112, 11).Value = "23/11/2023 05:12"
$ws.Cells.Item(112, 12).Value = 1.9
$ws.Cells.Item(112, 13).Value = "24/11/2023 16:51"
$ws.Cells.Item(112, 14).Value = 3.36
$ws.Cells.Item(112, 15).Value = "23/11/2023 05:12"
$ws.Cells.Item(112, 16).Value = 3.6
$ws.Cells.Item(112, 17).Value = "24/11/2023 16:51"
$ws.Cells.Item(112, 18).Value = 3.67
$ws.Cells.Item(112, 19).Value = "23/11/2023 05:12"
$ws.Cells.Item(112, 20).Value = 3.72
$ws.Cells.Item(112, 21).Value = "24/11/2023 16:51"
$ws.Cells.Item(112, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-spartak-subotica/xd8gwneP/"
